$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Countries")

# R1 values for Brobdingnag, Carpania, Dinotopia, Erewhon changed from 0 to 10
$ws.Range("C3").Value = 10
$ws.Range("C4").Value = 10
$ws.Range("C5").Value = 10
$ws.Range("C6").Value = 10

# Active selection moved to C7
$ws.Range("C7").Select() | Out-Null
